$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A1 with the modified text
$ws.Range("A1").Value = "Donnée A1 - modif 12h10"

# Reset the active selection to A1 (the default), clearing the previous C5 selection
$ws.Range("A1").Select()
